$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13, pushing the existing weekly records
# (previously in rows 13-54) down to rows 14-55.
$ws.Rows(13).Insert()

# Populate the new row 13 with the latest weekly price record.
$ws.Cells.Item(13, 1).Value = 4
$ws.Cells.Item(13, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(13, 3).Value = "Los Lagos"
$ws.Cells.Item(13, 4).Value = 44659
$ws.Cells.Item(13, 5).Value = 10
$ws.Cells.Item(13, 6).Value = 100112043
$ws.Cells.Item(13, 7).Value = "Pepino dulce"
$ws.Cells.Item(13, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 80
$ws.Cells.Item(13, 11).Value = 18000
$ws.Cells.Item(13, 12).Value = 18000
$ws.Cells.Item(13, 13).Value = 18000
$ws.Cells.Item(13, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(13, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(13, 16).Value = 1000
$ws.Cells.Item(13, 17).Value = 18
$ws.Cells.Item(13, 18).Value = "Hortaliza"
